$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values per row for columns B,C,E,F,G,H,I,K,L
# (columns D, J, M, N, O are unchanged and remain 0)
$data = @{
    2 = @{ "B"=12.9028481088006; "C"=6.328011044005836; "E"=11.13693434549475; "F"=16.86991607391233; "G"=41.85796577197106; "H"=17.64577311108699; "I"=27.69421552867239; "K"=10.48891585163557; "L"=9.985178365051512 }
    3 = @{ "B"=12.67711030025141; "C"=6.260596087588456; "E"=11.14247955847833; "F"=15.89584955866808; "G"=41.98770798483827; "H"=17.70758112158418; "I"=27.8033579276712; "K"=10.33935593532371; "L"=9.974354823361475 }
    4 = @{ "B"=12.53908135891781; "C"=6.218121665259912; "E"=11.14777588125694; "F"=15.26997757108491; "G"=42.07962878965829; "H"=17.74852778319295; "I"=27.87574070559022; "K"=10.24845791890527; "L"=9.969484579386549 }
    5 = @{ "B"=12.48305391754437; "C"=6.200547786865634; "E"=11.15041032109572; "F"=15.00819731993403; "G"=42.12015471211669; "H"=17.76596668253217; "I"=27.90658520498286; "K"=10.21169555269707; "L"=9.967948245778063 }
    6 = @{ "B"=12.47376615000796; "C"=6.197613838683204; "E"=11.150876537516; "F"=14.96433081551593; "G"=42.12706882222449; "H"=17.76890784732645; "I"=27.91178826041321; "K"=10.20560936551404; "L"=9.967720266348955 }
    7 = @{ "B"=12.53832475981577; "C"=6.217885722626796; "E"=11.14780948178412; "F"=15.26647399323137; "G"=42.08016293696591; "H"=17.7487599230835; "I"=27.87615123041209; "K"=10.24796093877493; "L"=9.969462042301796 }
    8 = @{ "B"=12.82493965357057; "C"=6.304994782557253; "E"=11.13845408989004; "F"=16.53996406344768; "G"=41.90014778259798; "H"=17.66646224829217; "I"=27.73073223506999; "K"=10.43718104349344; "L"=9.981079117116266 }
    9 = @{ "B"=13.38827038401731; "C"=6.46695018556598; "E"=11.13508956486825; "F"=19.00274580682531; "G"=41.64504582425904; "H"=17.52887957523238; "I"=27.4882608019614; "K"=10.81366607304106; "L"=10.01785450299068 }
    10 = @{ "B"=13.79864919051362; "C"=6.580167379505299; "E"=11.14170665723819; "F"=20.67494806633232; "G"=41.51813564390773; "H"=17.442347490364; "I"=27.33627481897953; "K"=11.09095589930993; "L"=10.05326448242085 }
    11 = @{ "B"=13.983719880657; "C"=6.63034541802653; "E"=11.14667785489218; "F"=21.3917225636224; "G"=41.47369923969604; "H"=17.40614978490266; "I"=27.27284064632861; "K"=11.21670944027825; "L"=10.07115995865406 }
    12 = @{ "B"=14.05350152718983; "C"=6.649149938433747; "E"=11.14884092583745; "F"=21.65686569030329; "G"=41.4587956362907; "H"=17.39289873077952; "I"=27.2496426065911; "K"=11.26422964310852; "L"=10.07819004956744 }
    13 = @{ "B"=14.03848722895055; "C"=6.64510889950379; "E"=11.14836261496812; "F"=21.60004134736742; "G"=41.46191967262286; "H"=17.39573227634449; "I"=27.25460206811976; "K"=11.25400045127625; "L"=10.07666478023551 }
    14 = @{ "B"=13.98946726497818; "C"=6.631896456864844; "E"=11.14685019737017; "F"=21.4136618050453; "G"=41.47243449890068; "H"=17.40505046375057; "I"=27.27091561771951; "K"=11.2206212072509; "L"=10.07173326997411 }
    15 = @{ "B"=13.95939999153409; "C"=6.623777640673932; "E"=11.14596029371182; "F"=21.29868154950795; "G"=41.47912595573858; "H"=17.41081756535221; "I"=27.28101541163378; "K"=11.20016115263169; "L"=10.06874547595689 }
    16 = @{ "B"=13.7865159858237; "C"=6.576860922339732; "E"=11.14142112253766; "F"=20.62722412089977; "G"=41.52130822604136; "H"=17.44477688723056; "I"=27.34053542400306; "K"=11.08272600423546; "L"=10.05213063873118 }
    17 = @{ "B"=13.67999427775943; "C"=6.547735130708328; "E"=11.13913781178201; "F"=20.20408069597325; "G"=41.5505996114868; "H"=17.46642143096306; "I"=27.37851215428022; "K"=11.01055227259254; "L"=10.04239330468089 }
    18 = @{ "B"=13.6185779236241; "C"=6.530858389777566; "E"=11.13800919050546; "F"=19.95656407809801; "G"=41.56869806556448; "H"=17.47916875868232; "I"=27.40089218382527; "K"=10.96900599141861; "L"=10.03696110633354 }
    19 = @{ "B"=13.59776014928153; "C"=6.525123051710594; "E"=11.13765881758257; "F"=19.87204792380568; "G"=41.57504035067227; "H"=17.48353593100309; "I"=27.40856179342916; "K"=10.95493469157663; "L"=10.03515089398457 }
    20 = @{ "B"=13.69134950916539; "C"=6.550848537468493; "E"=11.13936177022705; "F"=20.24955283636154; "G"=41.54735195046456; "H"=17.46408648792068; "I"=27.37441388800115; "K"=11.0182391192198; "L"=10.0434124505348 }
    21 = @{ "B"=14.0038743060147; "C"=6.635782661966223; "E"=11.14728682815244; "F"=21.46857628470577; "G"=41.46929374925077; "H"=17.40230109673092; "I"=27.26610157211765; "K"=11.23042855820201; "L"=10.07317492478567 }
    22 = @{ "B"=14.20634569228063; "C"=6.690142082582573; "E"=11.15410100607123; "F"=22.22866616901552; "G"=41.42949443431876; "H"=17.36458047687577; "I"=27.20011217007474; "K"=11.36850466022736; "L"=10.09410196281504 }
    23 = @{ "B"=14.0984677468149; "C"=6.661236673911643; "E"=11.15031508731551; "F"=21.82633154458858; "G"=41.44970624062795; "H"=17.38446900670654; "I"=27.23489194637663; "K"=11.29488016950244; "L"=10.08279903758068 }
    24 = @{ "B"=13.68621635074343; "C"=6.549441376778891; "E"=11.13925994511421; "F"=20.22900810905287; "G"=41.54881629884042; "H"=17.46514117003489; "I"=27.37626501139582; "K"=11.01476405642484; "L"=10.04295117759316 }
    25 = @{ "B"=13.23618770523022; "C"=6.424124588967123; "E"=11.13439978431858; "F"=18.34778573295695; "G"=41.70348721102989; "H"=17.5635470754515; "I"=27.54927162729305; "K"=10.71150905569843; "L"=10.00642096636884 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
